# TournRPG-301: バトル終了時にHP自動回復
# Add two new "auto" skills to the skill_auto sheet:
#   SKILL507 天の加護     (AUTO, win_hp = 15, 戦闘終了後にHPが15回復)
#   SKILL508 ソウルヒール (AUTO, win_mp = 10, 戦闘終了後にMPが10回復)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("skill_auto")

# Duplicate the formatting of the last existing data row (row 8) down into
# the two new rows so borders / fills / number formats stay consistent with
# the rest of the table instead of creating brand-new style records.
$ws.Range("A8:O8").Copy($ws.Range("A9:O9")) | Out-Null
$ws.Range("A8:O8").Copy($ws.Range("A10:O10")) | Out-Null

$ws.Rows(9).RowHeight = 20
$ws.Rows(10).RowHeight = 20

# Row 8's "rec_mp" sample value (column G) got copied into both new rows;
# neither new skill uses rec_mp, so clear it back out again.
$ws.Range("G9").ClearContents() | Out-Null
$ws.Range("G10").ClearContents() | Out-Null

# --- New skill: SKILL507 / 天の加護 (AUTO, win_hp = 15) ---
$ws.Range("A9").Value = "SKILL507"
$ws.Range("B9").Value = "天の加護"
$ws.Range("C9").Value = "AUTO"
$ws.Range("I9").NumberFormat = "0"
$ws.Range("I9").Value = 15
$ws.Range("O9").Value = "戦闘終了後にHPが15回復"

# --- New skill: SKILL508 / ソウルヒール (AUTO, win_mp = 10) ---
$ws.Range("A10").Value = "SKILL508"
$ws.Range("B10").Value = "ソウルヒール"
$ws.Range("C10").Value = "AUTO"
$ws.Range("J10").Value = 10
$ws.Range("O10").Value = "戦闘終了後にMPが10回復"
